# Updates loading_percent.xlsx results for "case with 380 kV" (Case_5_16)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 2).Value = 16.97360316259546
$ws.Cells.Item(2, 3).Value = 8.927199038821946
$ws.Cells.Item(2, 4).Value = 8.179484323685445
$ws.Cells.Item(2, 5).Value = 12.71991215803201
$ws.Cells.Item(2, 6).Value = 34.81600840799248
$ws.Cells.Item(2, 9).Value = 25.90020916364046
$ws.Cells.Item(2, 10).Value = 9.934216296803422
$ws.Cells.Item(2, 12).Value = 11.31073906751209
$ws.Cells.Item(2, 13).Value = 16.43699828793629
$ws.Cells.Item(2, 15).Value = 26.75469176677674

$ws.Cells.Item(3, 2).Value = 16.47134532421442
$ws.Cells.Item(3, 3).Value = 8.610807771909453
$ws.Cells.Item(3, 4).Value = 8.181699412178238
$ws.Cells.Item(3, 5).Value = 12.75325384499578
$ws.Cells.Item(3, 6).Value = 34.94969786146713
$ws.Cells.Item(3, 9).Value = 26.04598948743855
$ws.Cells.Item(3, 10).Value = 9.952701352754108
$ws.Cells.Item(3, 12).Value = 11.30793128581077
$ws.Cells.Item(3, 13).Value = 16.31609476631339
$ws.Cells.Item(3, 15).Value = 26.87018873949422

$ws.Cells.Item(4, 2).Value = 16.15609732352316
$ws.Cells.Item(4, 3).Value = 8.409388972004225
$ws.Cells.Item(4, 4).Value = 8.183768361102507
$ws.Cells.Item(4, 5).Value = 12.77484751703677
$ws.Cells.Item(4, 6).Value = 35.04026403197751
$ws.Cells.Item(4, 9).Value = 26.14085487560599
$ws.Cells.Item(4, 10).Value = 9.964656550594725
$ws.Cells.Item(4, 12).Value = 11.30733878506446
$ws.Cells.Item(4, 13).Value = 16.24291824750643
$ws.Cells.Item(4, 15).Value = 26.94731461301953

$ws.Cells.Item(5, 2).Value = 16.02608497996088
$ws.Cells.Item(5, 3).Value = 8.325597113346564
$ws.Cells.Item(5, 4).Value = 8.184790320582309
$ws.Cells.Item(5, 5).Value = 12.78392994963343
$ws.Cells.Item(5, 6).Value = 35.07929759272006
$ws.Cells.Item(5, 9).Value = 26.18086078127853
$ws.Cells.Item(5, 10).Value = 9.969681047318959
$ws.Cells.Item(5, 12).Value = 11.30738307239289
$ws.Cells.Item(5, 13).Value = 16.21338629803021
$ws.Cells.Item(5, 15).Value = 26.98030246486622

$ws.Cells.Item(6, 2).Value = 16.00440862042997
$ws.Cells.Item(6, 3).Value = 8.311582857480191
$ws.Cells.Item(6, 4).Value = 8.184970834923874
$ws.Cells.Item(6, 5).Value = 12.78545518576075
$ws.Cells.Item(6, 6).Value = 35.08590739916288
$ws.Cells.Item(6, 9).Value = 26.1875851400909
$ws.Cells.Item(6, 10).Value = 9.9705245944791
$ws.Cells.Item(6, 12).Value = 11.30740771814485
$ws.Cells.Item(6, 13).Value = 16.20850056857322
$ws.Cells.Item(6, 15).Value = 26.98587410431335

$ws.Cells.Item(7, 2).Value = 16.15434994644212
$ws.Cells.Item(7, 3).Value = 8.408265736152718
$ws.Cells.Item(7, 4).Value = 8.183781418689755
$ws.Cells.Item(7, 5).Value = 12.77496885968653
$ws.Cells.Item(7, 6).Value = 35.04078184681644
$ws.Cells.Item(7, 9).Value = 26.14138895226667
$ws.Cells.Item(7, 10).Value = 9.964723694017948
$ws.Cells.Item(7, 12).Value = 11.30733822384992
$ws.Cells.Item(7, 13).Value = 16.24251877439439
$ws.Cells.Item(7, 15).Value = 26.94775319294516

$ws.Cells.Item(8, 2).Value = 16.80195149638622
$ws.Cells.Item(8, 3).Value = 8.819641057639828
$ws.Cells.Item(8, 4).Value = 8.180101320800274
$ws.Cells.Item(8, 5).Value = 12.73117601561061
$ws.Cells.Item(8, 6).Value = 34.86034128753442
$ws.Cells.Item(8, 9).Value = 25.94936309808237
$ws.Cells.Item(8, 10).Value = 9.940464585190103
$ws.Cells.Item(8, 12).Value = 11.30953685048104
$ws.Cells.Item(8, 13).Value = 16.39510425348541
$ws.Cells.Item(8, 15).Value = 26.79322427583549

$ws.Cells.Item(9, 2).Value = 18.01029439511294
$ws.Cells.Item(9, 3).Value = 9.566254641836162
$ws.Cells.Item(9, 4).Value = 8.178482941893748
$ws.Cells.Item(9, 5).Value = 12.65416392544987
$ws.Cells.Item(9, 6).Value = 34.57401742846884
$ws.Cells.Item(9, 9).Value = 25.61526115507191
$ws.Cells.Item(9, 10).Value = 9.897674795464081
$ws.Cells.Item(9, 12).Value = 11.32277126833785
$ws.Cells.Item(9, 13).Value = 16.7017474894972
$ws.Cells.Item(9, 15).Value = 26.53961049476534

$ws.Cells.Item(10, 2).Value = 18.85185638328617
$ws.Cells.Item(10, 3).Value = 10.07430846079906
$ws.Cells.Item(10, 4).Value = 8.180671025018032
$ws.Cells.Item(10, 5).Value = 12.60293877286677
$ws.Cells.Item(10, 6).Value = 34.40512485986806
$ws.Cells.Item(10, 9).Value = 25.39563289959245
$ws.Cells.Item(10, 10).Value = 9.869123757873982
$ws.Cells.Item(10, 12).Value = 11.33785814413255
$ws.Cells.Item(10, 13).Value = 16.93023627193102
$ws.Cells.Item(10, 15).Value = 26.38359764035937

$ws.Cells.Item(11, 2).Value = 19.2231080856531
$ws.Cells.Item(11, 3).Value = 10.29599013201312
$ws.Cells.Item(11, 4).Value = 8.18239165465388
$ws.Cells.Item(11, 5).Value = 12.58078763064008
$ws.Cells.Item(11, 6).Value = 34.3373581756747
$ws.Cells.Item(11, 9).Value = 25.30131964981429
$ws.Cells.Item(11, 10).Value = 9.856755844733254
$ws.Cells.Item(11, 12).Value = 11.34586855992464
$ws.Cells.Item(11, 13).Value = 17.03458952136335
$ws.Cells.Item(11, 15).Value = 26.31924689785494

$ws.Cells.Item(12, 2).Value = 19.36191317263294
$ws.Cells.Item(12, 3).Value = 10.37853314289665
$ws.Cells.Item(12, 4).Value = 8.183146759054214
$ws.Cells.Item(12, 5).Value = 12.57256436227229
$ws.Cells.Item(12, 6).Value = 34.31300482435212
$ws.Cells.Item(12, 9).Value = 25.26641014116415
$ws.Cells.Item(12, 10).Value = 9.85216114740513
$ws.Cells.Item(12, 12).Value = 11.34906522899129
$ws.Cells.Item(12, 13).Value = 17.07414040991739
$ws.Cells.Item(12, 15).Value = 26.29583428928676

$ws.Cells.Item(13, 2).Value = 19.3320998748416
$ws.Cells.Item(13, 3).Value = 10.36081910070703
$ws.Cells.Item(13, 4).Value = 8.182979541154891
$ws.Cells.Item(13, 5).Value = 12.57432806765871
$ws.Cells.Item(13, 6).Value = 34.31819148303419
$ws.Cells.Item(13, 9).Value = 25.27389272952699
$ws.Cells.Item(13, 10).Value = 9.85314675635545
$ws.Cells.Item(13, 12).Value = 11.34836953401607
$ws.Cells.Item(13, 13).Value = 17.06562127386701
$ws.Cells.Item(13, 15).Value = 26.30083405467176

$ws.Cells.Item(14, 2).Value = 19.23456395004566
$ws.Cells.Item(14, 3).Value = 10.30280935732223
$ws.Cells.Item(14, 4).Value = 8.182451706874829
$ws.Cells.Item(14, 5).Value = 12.58010779670304
$ws.Cells.Item(14, 6).Value = 34.33532835965293
$ws.Cells.Item(14, 9).Value = 25.29843149198666
$ws.Cells.Item(14, 10).Value = 9.856376059590247
$ws.Cells.Item(14, 12).Value = 11.34612828899236
$ws.Cells.Item(14, 13).Value = 17.0378428452954
$ws.Cells.Item(14, 15).Value = 26.31730155577251

$ws.Cells.Item(15, 2).Value = 19.17458538522649
$ws.Cells.Item(15, 3).Value = 10.2670927299747
$ws.Cells.Item(15, 4).Value = 8.182141854202431
$ws.Cells.Item(15, 5).Value = 12.5836695007009
$ws.Cells.Item(15, 6).Value = 34.34599573417454
$ws.Cells.Item(15, 9).Value = 25.31356701417605
$ws.Cells.Item(15, 10).Value = 9.858365648632519
$ws.Cells.Item(15, 12).Value = 11.34477667666282
$ws.Cells.Item(15, 13).Value = 17.02083152533726
$ws.Cells.Item(15, 15).Value = 26.32751293265006

$ws.Cells.Item(16, 2).Value = 18.82735111915031
$ws.Cells.Item(16, 3).Value = 10.05962708846323
$ws.Cells.Item(16, 4).Value = 8.180573108257972
$ws.Cells.Item(16, 5).Value = 12.60440951310378
$ws.Cells.Item(16, 6).Value = 34.40973639354085
$ws.Cells.Item(16, 9).Value = 25.4019091167724
$ws.Cells.Item(16, 10).Value = 9.869944473206846
$ws.Cells.Item(16, 12).Value = 11.33735758938083
$ws.Cells.Item(16, 13).Value = 16.92342272677941
$ws.Cells.Item(16, 15).Value = 26.38793664144491

$ws.Cells.Item(17, 2).Value = 18.61128290216377
$ws.Cells.Item(17, 3).Value = 9.929902215621984
$ws.Cells.Item(17, 4).Value = 8.179795963706656
$ws.Cells.Item(17, 5).Value = 12.61742725882003
$ws.Cells.Item(17, 6).Value = 34.45116429197938
$ws.Cells.Item(17, 9).Value = 25.45753761538349
$ws.Cells.Item(17, 10).Value = 9.877206238496395
$ws.Cells.Item(17, 12).Value = 11.33309890169621
$ws.Cells.Item(17, 13).Value = 16.86375362990942
$ws.Cells.Item(17, 15).Value = 26.42670284410597

$ws.Cells.Item(18, 2).Value = 18.48592237212871
$ws.Cells.Item(18, 3).Value = 9.854402220634402
$ws.Cells.Item(18, 4).Value = 8.179417282403817
$ws.Cells.Item(18, 5).Value = 12.62502314116985
$ws.Cells.Item(18, 6).Value = 34.47584532152869
$ws.Cells.Item(18, 9).Value = 25.49006043894931
$ws.Cells.Item(18, 10).Value = 9.881441404329186
$ws.Cells.Item(18, 12).Value = 11.33075751410533
$ws.Cells.Item(18, 13).Value = 16.82947385590326
$ws.Cells.Item(18, 15).Value = 26.44962305711368

$ws.Cells.Item(19, 2).Value = 18.44329500360873
$ws.Cells.Item(19, 3).Value = 9.82868861965944
$ws.Cells.Item(19, 4).Value = 8.179300822464274
$ws.Cells.Item(19, 5).Value = 12.62761361857079
$ws.Cells.Item(19, 6).Value = 34.48434820316577
$ws.Cells.Item(19, 9).Value = 25.50116258304424
$ws.Cells.Item(19, 10).Value = 9.882885399779981
$ws.Cells.Item(19, 12).Value = 11.3299833763185
$ws.Cells.Item(19, 13).Value = 16.81787499185869
$ws.Cells.Item(19, 15).Value = 26.45749032282528

$ws.Cells.Item(20, 2).Value = 18.63439677267416
$ws.Cells.Item(20, 3).Value = 9.943803669946563
$ws.Cells.Item(20, 4).Value = 8.179871627432547
$ws.Cells.Item(20, 5).Value = 12.61603028156928
$ws.Cells.Item(20, 6).Value = 34.44666592221559
$ws.Cells.Item(20, 9).Value = 25.45156134650782
$ws.Cells.Item(20, 10).Value = 9.876427170939708
$ws.Cells.Item(20, 12).Value = 11.33354107060529
$ws.Cells.Item(20, 13).Value = 16.87010151438104
$ws.Cells.Item(20, 15).Value = 26.42251161942837

$ws.Cells.Item(21, 2).Value = 19.26326177629152
$ws.Cells.Item(21, 3).Value = 10.31988664142107
$ws.Cells.Item(21, 4).Value = 8.182603940706846
$ws.Cells.Item(21, 5).Value = 12.57840568059617
$ws.Cells.Item(21, 6).Value = 34.33025929383965
$ws.Cells.Item(21, 9).Value = 25.29120201764229
$ws.Cells.Item(21, 10).Value = 9.8554251291678
$ws.Cells.Item(21, 12).Value = 11.34678217892014
$ws.Cells.Item(21, 13).Value = 17.04600130159353
$ws.Cells.Item(21, 15).Value = 26.31243868755584

$ws.Cells.Item(22, 2).Value = 19.66384199829099
$ws.Cells.Item(22, 3).Value = 10.55748345235322
$ws.Cells.Item(22, 4).Value = 8.184992802608589
$ws.Cells.Item(22, 5).Value = 12.55477661080576
$ws.Cells.Item(22, 6).Value = 34.26180951221474
$ws.Cells.Item(22, 9).Value = 25.19108943265691
$ws.Cells.Item(22, 10).Value = 9.842216273473266
$ws.Cells.Item(22, 12).Value = 11.35638700127522
$ws.Cells.Item(22, 13).Value = 17.16115211202996
$ws.Cells.Item(22, 15).Value = 26.24607140104283

$ws.Cells.Item(23, 2).Value = 19.4510330340646
$ws.Cells.Item(23, 3).Value = 10.43143715227782
$ws.Cells.Item(23, 4).Value = 8.183662890796073
$ws.Cells.Item(23, 5).Value = 12.56730020367136
$ws.Cells.Item(23, 6).Value = 34.29764279745001
$ws.Cells.Item(23, 9).Value = 25.24409208464485
$ws.Cells.Item(23, 10).Value = 9.849218898403031
$ws.Cells.Item(23, 12).Value = 11.35117428738685
$ws.Cells.Item(23, 13).Value = 17.09968471494032
$ws.Cells.Item(23, 15).Value = 26.2809818704443

$ws.Cells.Item(24, 2).Value = 18.62395053392398
$ws.Cells.Item(24, 3).Value = 9.937521682752342
$ws.Cells.Item(24, 4).Value = 8.179837207664225
$ws.Cells.Item(24, 5).Value = 12.61666150666592
$ws.Cells.Item(24, 6).Value = 34.44869694556161
$ws.Cells.Item(24, 9).Value = 25.45426153165137
$ws.Cells.Item(24, 10).Value = 9.876779199579868
$ws.Cells.Item(24, 12).Value = 11.3333408326582
$ws.Cells.Item(24, 13).Value = 16.86723155671618
$ws.Cells.Item(24, 15).Value = 26.42440450057802

$ws.Cells.Item(25, 2).Value = 17.69094625406601
$ws.Cells.Item(25, 3).Value = 9.371138443177534
$ws.Cells.Item(25, 4).Value = 8.178325005762206
$ws.Cells.Item(25, 5).Value = 12.67405371194172
$ws.Cells.Item(25, 6).Value = 34.64421464938693
$ws.Cells.Item(25, 9).Value = 25.70110424617187
$ws.Cells.Item(25, 10).Value = 9.908741582628512
$ws.Cells.Item(25, 12).Value = 11.31824330138003
$ws.Cells.Item(25, 13).Value = 16.61813367952537
$ws.Cells.Item(25, 15).Value = 26.6029094372265

Write-Host "Updated 240 cells (rows 2-25, cols B:F,I:J,L:M,O) for case with 380 kV"